# Update the extracted_profiles worksheet:
# - Insert a new profile (Madison McGuire) at row 5, which pushes the
#   previously-listed Savanah Nunes / Steve Greene / Amy Calabretta / Angela Harley
#   rows down by one.
# - Replace the former Brandi Chastain / Brian Kelly / Carly Schoepflin rows
#   (rows 9-11) with Angela Harley (shifted down), Ann Barrington and Ashley Baker.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: Madison McGuire
$ws.Range("A5").Value2 = "https://www.uidaho.edu/cals/people/madison-mcguire"
$ws.Range("C5").Value2 = "Madison"
$ws.Range("E5").Value2 = "McGuire"
$ws.Range("F5").Value2 = "Administrative Specialist"

# Row 6: Savanah Nunes
$ws.Range("A6").Value2 = "https://www.uidaho.edu/cals/people/savanah-nunes"
$ws.Range("C6").Value2 = "Savanah"
$ws.Range("E6").Value2 = "Nunes"
$ws.Range("F6").Value2 = "Media and Communications Manager"

# Row 7: Steve Greene
$ws.Range("A7").Value2 = "https://www.uidaho.edu/cals/people/steve-greene"
$ws.Range("C7").Value2 = "Steve"
$ws.Range("E7").Value2 = "Greene"
$ws.Range("F7").Value2 = "Program Manager"

# Row 8: Amy Calabretta
$ws.Range("A8").Value2 = "https://www.uidaho.edu/cals/people/amy-calabretta"
$ws.Range("C8").Value2 = "Amy"
$ws.Range("E8").Value2 = "Calabretta"
$ws.Range("F8").Value2 = "Interim Director of Communications & Strategic Initiatives"
$ws.Range("J8").ClearContents()

# Row 9: Angela Harley
$ws.Range("A9").Value2 = "https://www.uidaho.edu/cals/people/angela-harley"
$ws.Range("C9").Value2 = "Angela"
$ws.Range("E9").Value2 = "Harley"
$ws.Range("F9").ClearContents()
$ws.Range("J9").Value2 = "College of Agricultural & Life Sciences"

# Row 10: Ann Barrington (replaces Brian Kelly)
$ws.Range("A10").Value2 = "https://www.uidaho.edu/cals/people/ann-barrington"
$ws.Range("C10").Value2 = "Ann"
$ws.Range("E10").Value2 = "Barrington"
$ws.Range("F10").Value2 = "Senior Director of Development"

# Row 11: Ashley Baker (replaces Carly Schoepflin)
$ws.Range("A11").Value2 = "https://www.uidaho.edu/cals/people/ashley-baker"
$ws.Range("C11").Value2 = "Ashley"
$ws.Range("E11").Value2 = "Baker"
$ws.Range("F11").Value2 = "Instructional Media Designer"
